# Update "want to go" counts (column F) on sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1022
$ws1.Range("F4").Value = 13454
$ws1.Range("F6").Value = 1016
$ws1.Range("F7").Value = 10
$ws1.Range("F8").Value = 1732
$ws1.Range("F9").Value = 130
$ws1.Range("F10").Value = 120
$ws1.Range("F11").Value = 77
$ws1.Range("F13").Value = 31
$ws1.Range("F14").Value = 13442
$ws1.Range("F15").Value = 334
$ws1.Range("F16").Value = 593
$ws1.Range("F17").Value = 8941
$ws1.Range("F19").Value = 8007
$ws1.Range("F20").Value = 247
$ws1.Range("F22").Value = 142
$ws1.Range("F27").Value = 1019
$ws1.Range("F32").Value = 169
$ws1.Range("F34").Value = 92

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1022
$ws4.Range("F5").Value = 13454
$ws4.Range("F7").Value = 1016
$ws4.Range("F8").Value = 10
$ws4.Range("F9").Value = 1732
$ws4.Range("F10").Value = 130
$ws4.Range("F11").Value = 120
$ws4.Range("F12").Value = 77
$ws4.Range("F14").Value = 31
$ws4.Range("F15").Value = 13442
$ws4.Range("F16").Value = 334
$ws4.Range("F17").Value = 593
$ws4.Range("F18").Value = 8941
$ws4.Range("F20").Value = 8007
$ws4.Range("F21").Value = 247
$ws4.Range("F23").Value = 142
$ws4.Range("F28").Value = 1019
$ws4.Range("F35").Value = 169
$ws4.Range("F37").Value = 92
